$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 2587.3333
$ws.Range("J88").Value = 2596
$ws.Range("L88").Value = 2596
$ws.Range("N88").Value = -3408
$ws.Range("H91").Value = 2587.3333
$ws.Range("J91").Value = 2596
$ws.Range("L91").Value = 2596
$ws.Range("N91").Value = -5404
$ws.Range("H106").Value = 20799.445
$ws.Range("I106").Value = 16249.833
$ws.Range("K106").Value = 16249.833
$ws.Range("M106").Value = -15618.833
$ws.Range("H125").Value = 747.5
$ws.Range("I125").Value = 497
$ws.Range("K125").Value = 4473
$ws.Range("M125").Value = -2013
$ws.Range("H129").Value = 2828.5386
$ws.Range("I129").Value = 1855.4
$ws.Range("K129").Value = 5566.200000000001
$ws.Range("M129").Value = -566.2000000000007
$ws.Range("H132").Value = 5369.7085
$ws.Range("I132").Value = 5721.9414
$ws.Range("J132").Value = 4514.2856
$ws.Range("K132").Value = 17165.8242
$ws.Range("L132").Value = 13542.8568
$ws.Range("M132").Value = -14635.8242
$ws.Range("N132").Value = -18602.8568
$ws.Range("H135").Value = 620
$ws.Range("I135").Value = 620
$ws.Range("K135").Value = 5580
$ws.Range("M135").Value = -3045
$ws.Range("H138").Value = 1703.4286
$ws.Range("J138").Value = 2795
$ws.Range("L138").Value = 8385
$ws.Range("N138").Value = -18665
$ws.Range("H141").Value = 14899.25
$ws.Range("I141").Value = 19998.5
$ws.Range("K141").Value = 59995.5
$ws.Range("M141").Value = -54815.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2670.889
$ws.Range("I97").Value = 1754.75
$ws.Range("K97").Value = 1754.75
$ws.Range("M97").Value = -1258.75
$ws.Range("H131").Value = 57958
$ws.Range("J131").Value = 57958
$ws.Range("L131").Value = 57958
$ws.Range("N131").Value = -68038

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1654.0476
$ws.Range("I94").Value = 4161.6665
$ws.Range("K94").Value = 4161.6665
$ws.Range("M94").Value = -3710.6665
$ws.Range("H134").Value = 1562.5238
$ws.Range("I134").Value = 1562.5238
$ws.Range("K134").Value = 4687.5714
$ws.Range("M134").Value = -2152.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1391
$ws.Range("I16").Value = 1391
$ws.Range("K16").Value = 1391
$ws.Range("M16").Value = -1104
$ws.Range("H31").Value = 3276.6667
$ws.Range("I31").Value = 2358
$ws.Range("K31").Value = 2358
$ws.Range("M31").Value = -2063
$ws.Range("H34").Value = 3276.6667
$ws.Range("I34").Value = 2358
$ws.Range("K34").Value = 2358
$ws.Range("M34").Value = -2156
$ws.Range("H58").Value = 2189.3
$ws.Range("I58").Value = 1358.6
$ws.Range("K58").Value = 1358.6
$ws.Range("M58").Value = -1155.6
$ws.Range("H99").Value = 6082.5
$ws.Range("I99").Value = 6082.5
$ws.Range("K99").Value = 6082.5
$ws.Range("M99").Value = -4584.5
$ws.Range("H113").Value = 1391
$ws.Range("I113").Value = 1391
$ws.Range("K113").Value = 1391
$ws.Range("M113").Value = 779
$ws.Range("H126").Value = 6082.5
$ws.Range("I126").Value = 6082.5
$ws.Range("K126").Value = 18247.5
$ws.Range("M126").Value = -15777.5
$ws.Range("H132").Value = 3603.7856
$ws.Range("I132").Value = 3541.182
$ws.Range("J132").Value = 3833.3333
$ws.Range("K132").Value = 10623.546
$ws.Range("L132").Value = 11499.9999
$ws.Range("M132").Value = -8093.545999999998
$ws.Range("N132").Value = -16559.9999
$ws.Range("H136").Value = 2189.3
$ws.Range("I136").Value = 1358.6
$ws.Range("K136").Value = 4075.8
$ws.Range("M136").Value = -1525.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 330
$ws.Range("I20").Value = 330
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 990
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -763
$ws.Range("N20").ClearContents()
$ws.Range("H21").Value = 250
$ws.Range("I21").Value = 250
$ws.Range("K21").Value = 750
$ws.Range("M21").Value = -577
$ws.Range("H22").Value = 201
$ws.Range("I22").Value = 201
$ws.Range("K22").Value = 603
$ws.Range("M22").Value = -434
$ws.Range("H26").Value = 668.38464
$ws.Range("I26").Value = 410.1111
$ws.Range("J26").Value = 1249.5
$ws.Range("K26").Value = 1230.3333
$ws.Range("L26").Value = 3748.5
$ws.Range("M26").Value = -942.3333
$ws.Range("N26").Value = -4324.5
$ws.Range("H27").Value = 201
$ws.Range("I27").Value = 201
$ws.Range("K27").Value = 603
$ws.Range("M27").Value = -501
$ws.Range("H107").Value = 180
$ws.Range("J107").Value = 180
$ws.Range("L107").Value = 540
$ws.Range("N107").Value = -4380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3057.4
$ws.Range("I80").Value = 3221.75
$ws.Range("K80").Value = 3221.75
$ws.Range("M80").Value = -2223.75
$ws.Range("H83").Value = 3057.4
$ws.Range("I83").Value = 3221.75
$ws.Range("K83").Value = 16108.75
$ws.Range("M83").Value = -11116.75
$ws.Range("H97").Value = 1778.6666
$ws.Range("I97").Value = 1706.1538
$ws.Range("J97").Value = 2250
$ws.Range("K97").Value = 1706.1538
$ws.Range("L97").Value = 2250
$ws.Range("M97").Value = -1210.1538
$ws.Range("N97").Value = -3242
$ws.Range("H104").Value = 59399.6
$ws.Range("J104").Value = 65749.5
$ws.Range("L104").Value = 65749.5
$ws.Range("N104").Value = -72737.5
$ws.Range("H128").Value = 91000
$ws.Range("J128").Value = 91000
$ws.Range("L128").Value = 91000
$ws.Range("N128").Value = -100960
$ws.Range("H132").Value = 3117.3572
$ws.Range("I132").Value = 2998.25
$ws.Range("J132").Value = 3832
$ws.Range("K132").Value = 8994.75
$ws.Range("L132").Value = 11496
$ws.Range("M132").Value = -6464.75
$ws.Range("N132").Value = -16556

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1386.75
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H93").Value = 3473.75
$ws.Range("I93").Value = 3473.75
$ws.Range("K93").Value = 3473.75
$ws.Range("M93").Value = -2225.75
$ws.Range("H100").Value = 1504.2858
$ws.Range("I100").Value = 1504.2858
$ws.Range("K100").Value = 1504.2858
$ws.Range("M100").Value = -963.2858000000001
$ws.Range("H113").Value = 1386.75
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 62991.332
$ws.Range("J130").Value = 62991.332
$ws.Range("L130").Value = 62991.332
$ws.Range("N130").Value = -73031.33199999999
$ws.Range("H136").Value = 11166.333
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 15499.5
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 46498.5
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -51598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H122").Value = 3507.6
$ws.Range("I122").Value = 3341.7778
$ws.Range("K122").Value = 10025.3334
$ws.Range("M122").Value = -7575.3334
$ws.Range("H126").Value = 4482.4346
$ws.Range("I126").Value = 4240.5625
$ws.Range("K126").Value = 12721.6875
$ws.Range("M126").Value = -10251.6875
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H136").Value = 14181.941
$ws.Range("I136").Value = 8886.267
$ws.Range("K136").Value = 26658.801
$ws.Range("M136").Value = -24108.801
